$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update "new Madigan bike hours" — Riders/Average values for 28, 29, 31 Oct 2016
$ws.Range("C2").Value = 281
$ws.Range("D2").Value = 281

$ws.Range("C3").Value = 127
$ws.Range("D3").Value = 127

$ws.Range("C5").Value = 269
$ws.Range("D5").Value = 269
